$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): update F4, F5, F10
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1557
$ws1.Range("F5").Value = 2
$ws1.Range("F10").Value = 398

# Sheet "全部类型" (sheet4.xml): update F4, F5, F10
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1557
$ws4.Range("F5").Value = 2
$ws4.Range("F10").Value = 398
